$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "POST"
$ws.Range("C5").Value = "/proveedores/crear"
$ws.Range("D5").Value = "Crea un proveedor"

$ws.Range("B6").Value = "GET"
$ws.Range("C6").Value = "/productos/buscar/{subcadena}"
$ws.Range("D6").Value = "Recupera todos los productos cuyo nombre contenga la subcadena"

$ws.Range("B2:D6").Select() | Out-Null
